$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "Sign and date the TODI in front of 2 witnesses and a
# notary public." -> "{{multiple_must_sign}} and date the TODI in
# front of 2 witnesses and a notary public." (still bold, same as
# before) plus a bookmark ("_GoBack") marking the last edit point,
# right after the closing "}}".
# ------------------------------------------------------------------

$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Sign and date the TODI")) {
        $target1 = $p.Range
        break
    }
}

$paraStart = $target1.Start

# Replace the word "Sign" with the merge-field style tag text.
$signRange = $d.Range($paraStart, $paraStart + 4)
$signRange.Text = "{{multiple_must_sign}}"

# Work out where "}}" ends so we can drop the _GoBack bookmark there,
# matching where Word leaves it after the last text entry.
$afterTag = $paraStart + 22

$goBack = $d.Range($afterTag, $afterTag)
$d.Bookmarks.Add("_GoBack", $goBack)

# ------------------------------------------------------------------
# Change 2: "...Transfer Tax Stamp" box until later." -> "...until
# laters." (an extra run containing just "s" tacked on at the end of
# the paragraph).
# ------------------------------------------------------------------

$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.IndexOf("Transfer Tax Stamp") -ge 0) {
        $target2 = $p.Range
        break
    }
}

$insertPoint = $target2.End - 1
$tail = $d.Range($insertPoint, $insertPoint)
$tail.InsertAfter("s")
